# Apply the "add 2022-Q4 data" edit:
#  1. Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q3"),
#     cloned from the "2022-Q3" sheet so it inherits the same layout/styles.
#  2. Overwrite its data rows with the 2022-Q4 fund holdings.
#  3. Insert a new row for "2022-Q4" at the top of the data in the "总计"
#     summary sheet, shifting the existing rows down and renumbering the
#     index column.
#  4. Restore the originally-active sheet ("2020-Q4") as the selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone "2022-Q3" (currently sheet index 2) into a new sheet
# placed immediately before it, then rename it.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($q3Sheet, $null)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The source "2022-Q3" sheet only had 9 data rows (rows 2-10); 2022-Q4 has
# 16, so extend column A's index-cell style ("s=2") down through row 17
# before writing data, matching the style already present on rows 2-10.
$q4Sheet.Cells.Item(2, 1).Copy($q4Sheet.Range("A11:A17"))

# ---------------------------------------------------------------------
# Step 2: replace the data rows (2..) with the 2022-Q4 holdings.
# Columns: A=index, B=code, C=name, D=scale, E=stock position,
#          F=position pct, G=market value (billion), H=position rank
# Row1 headers were already copied correctly from "2022-Q3".
# ---------------------------------------------------------------------
$q4Data = @(
    @("163409", "兴全绿色投资混合（LOF）", "50.87", "91.54", "2.34", "1.1904", 10),
    @("070021", "嘉实主题新动力混合", "17.08", "93.12", "5.01", "0.8557", 1),
    @("000985", "嘉实逆向策略股票", "11.36", "93.39", "4.99", "0.5669", 1),
    @("012466", "嘉实策略精选混合A", "9.59", "93.89", "4.96", "0.4757", 1),
    @("001907", "国投瑞银境煊灵活配置混合A", "4.97", "90.03", "7.57", "0.3762", 7),
    @("001908", "国投瑞银境煊灵活配置混合C", "4.90", "90.03", "7.57", "0.3709", 7),
    @("010425", "国投瑞银开放视角精选混合A", "4.55", "91.41", "7.77", "0.3535", 6),
    @("010673", "兴全中证800六个月持有期指数增强A", "12.63", "95.84", "2.49", "0.3145", 10),
    @("015309", "国投瑞银境煊灵活配置混合E", "3.40", "90.03", "7.57", "0.2574", 7),
    @("013627", "华夏周期驱动混合C", "4.87", "86.45", "3.28", "0.1597", 9),
    @("013626", "华夏周期驱动混合A", "4.12", "86.45", "3.28", "0.1351", 9),
    @("014307", "嘉实多元动力混合A", "1.56", "93.22", "4.96", "0.0774", 2),
    @("012467", "嘉实策略精选混合C", "1.01", "93.89", "4.96", "0.0501", 1),
    @("010426", "国投瑞银开放视角精选混合C", "0.61", "91.41", "7.77", "0.0474", 6),
    @("010674", "兴全中证800六个月持有期指数增强C", "1.32", "95.84", "2.49", "0.0329", 10),
    @("014308", "嘉实多元动力混合C", "0.33", "93.22", "4.96", "0.0164", 2)
)

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r, 1).Value = $r - 2
    $q4Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4Sheet.Cells.Item($r, 2).ClearFormats()
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4Sheet.Cells.Item($r, 4).ClearFormats()
    $q4Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4Sheet.Cells.Item($r, 5).ClearFormats()
    $q4Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $q4Sheet.Cells.Item($r, 6).ClearFormats()
    $q4Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $q4Sheet.Cells.Item($r, 7).ClearFormats()
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet (sheet index 1): insert the new
# 2022-Q4 row at the top of the data and shift the rest down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()
# Match the index-column style ("s=2", same as the other A-column cells)
# used throughout the rest of the sheet.
$summary.Cells.Item(3, 1).Copy($summary.Cells.Item(2, 1))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 16
$summary.Cells.Item(2, 4).Value = 5.28

# renumber the index column (A) for all the rows that shifted down
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(9, 1).Value = 7
$summary.Cells.Item(10, 1).Value = 8

# ---------------------------------------------------------------------
# Step 4: restore "2020-Q4" (now the last sheet) as the active/selected
# tab, matching the workbook's original selection state.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
